# Update the "Producto" sheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Producto")

$ws.Range("B2").Value = "NuevoTestdfgfdg"
$ws.Range("B3").Value = "sdfsdfdsfsd"
$ws.Range("B4").Value = "carozzigterg"

# CONTENIDO's new value ("500") looks numeric; force it to stay a text
# cell (matching the original "50" text cell) by pre-formatting as Text.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "500"

$ws.Range("B6").Value = "GR"

# PRECIO's new value ("500") also looks numeric; same text-forcing trick.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "500"

$ws.Range("B9").Value = "Gaseosa"

# Remove the "Inventarios" sheet entirely
$inv = $wb.Worksheets.Item("Inventarios")
[void]$inv.Delete()
